$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.508.71'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.831.80'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '312.77'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4290'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = '0.3664'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.07295'
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").Value = '0.8672'
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("D11").Value = '20.70'
$ws.Range("D12").Value = '1.795.46'
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("D13").Value = '5.414'
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").Value = '6.529'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '0.06939'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '80.52'
$ws.Range("D18").Value = '0.000008932'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '15.44'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '27.517.00'
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("D22").Value = '5.139'
$ws.Range("E22").Value = '  +3.22%  '
$ws.Range("D23").Value = '10.86'
$ws.Range("E23").Value = '  +4.91%  '
$ws.Range("D24").Value = '2.038.54'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '154.54'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").Value = '19.01'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").Value = '5.146'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("D30").Value = '1.840'
$ws.Range("E30").Value = '  -2.12%  '
$ws.Range("D31").Value = '0.08875'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").Value = '0.7553'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '2.981'
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("D34").Value = '4.549'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '1.139'
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '1.091'
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").Value = '0.05333'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("D39").Value = '0.01941'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5104'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1667'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").Value = '6.589'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '8.376'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = '10.50'
$ws.Range("E45").Value = '  +0.92%  '
$ws.Range("D46").Value = '106.28'
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("D47").Value = '0.06504'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = '0.4693'
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = '1.619'
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("D51").Value = '64.00'
$ws.Range("E51").Value = '  -0.82%  '
